$d = $word.ActiveDocument

# The cell currently holds "{{bo_letter_names}}{{bo_correspond_add_vert}}"
# as one run of text. Replace it with the two placeholders split onto
# separate lines via a manual line break (^l), matching the diff which
# turns the trailing "}}{{" boundary into "</w:t><w:br/><w:t>".
$d.Content.Find.Execute(
    "{{bo_letter_names}}{{bo_correspond_add_vert}}",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "{{bo_letter_names}}^l{{bo_correspond_add_vert}}", 2) | Out-Null
